# Update the CDA Logical Model metadata sheet (ST.r2b related regeneration):
#  - bump Version string
#  - bump publication Date
#  - insert a new "Jurisdiction" property row (empty value) right after "Contact"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3, column B)
$ws.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8, column B)
$ws.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"

# Insert a new row for "Jurisdiction" after the "Contact" row (row 10), pushing
# "Description" and everything below down by one row.
$ws.Rows.Item(11).Insert()
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Match the formatting of the surrounding data rows (the freshly inserted row
# otherwise has no style applied).
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0
